# Updates the cryptos list: Price (column D) and Volume(1h) (column E)
# for rows 2-51, reflecting the latest scrape from the GitHub Actions job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cells whose new text is unambiguous (Excel will not reinterpret it as a number) ---
$ws.Range("D2").Value = "22.387.75"
$ws.Range("E2").Value = "  -4.51%  "
$ws.Range("D3").Value = "1.567.42"
$ws.Range("E3").Value = "  -4.74%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("E5").Value = "  +0.03%  "
$ws.Range("E6").Value = "  -2.96%  "
$ws.Range("E7").Value = "  -3.00%  "
$ws.Range("E8").Value = "  -1.19%  "
$ws.Range("E9").Value = "  -3.61%  "
$ws.Range("E10").Value = "  -3.65%  "
$ws.Range("E11").Value = "  -5.81%  "
$ws.Range("E12").Value = "  +0.05%  "
$ws.Range("E13").Value = "  -4.14%  "
$ws.Range("E14").Value = "  -5.02%  "
$ws.Range("E15").Value = "  -5.85%  "
$ws.Range("D16").Value = "1.571.06"
$ws.Range("E16").Value = "  -4.47%  "
$ws.Range("E17").Value = "  -5.62%  "
$ws.Range("E18").Value = "  -7.66%  "
$ws.Range("E19").Value = "  -3.28%  "
$ws.Range("E20").Value = "  +0.08%  "
$ws.Range("E21").Value = "  -7.14%  "
$ws.Range("E22").Value = "  -6.71%  "
$ws.Range("E23").Value = "  -5.02%  "
$ws.Range("E24").Value = "  -2.71%  "
$ws.Range("D25").Value = "22.393.75"
$ws.Range("E25").Value = "  -4.57%  "
$ws.Range("E26").Value = "  -3.48%  "
$ws.Range("E28").Value = "  -4.57%  "
$ws.Range("E29").Value = "  -4.94%  "
$ws.Range("E30").Value = "  -4.40%  "
$ws.Range("E31").Value = "  -5.59%  "
$ws.Range("D32").Value = "1.750.37"
$ws.Range("E32").Value = "  -4.18%  "
$ws.Range("E33").Value = "  +4.79%  "
$ws.Range("E34").Value = "  -9.36%  "
$ws.Range("E35").Value = "  -6.21%  "
$ws.Range("E36").Value = "  -9.83%  "
$ws.Range("E37").Value = "  -3.27%  "
$ws.Range("E38").Value = "  -5.99%  "
$ws.Range("E39").Value = "  -4.30%  "
$ws.Range("E40").Value = "  -3.63%  "
$ws.Range("E41").Value = "  -6.42%  "
$ws.Range("E42").Value = "  -8.45%  "
$ws.Range("E43").Value = "  -3.87%  "
$ws.Range("E44").Value = "  -7.33%  "
$ws.Range("E45").Value = "  -8.27%  "
$ws.Range("E46").Value = "  +0.02%  "
$ws.Range("E47").Value = "  -5.30%  "
$ws.Range("E48").Value = "  -3.25%  "
$ws.Range("E49").Value = "  -5.46%  "
$ws.Range("E50").Value = "  +8.29%  "
$ws.Range("E51").Value = "  -2.91%  "

# --- Price cells whose new text looks like a plain number (single "."): Excel would
# normally reinterpret these and drop formatting (e.g. "1.170" -> 1.17), so force them
# to Text, assign, then clear the transient formatting back to the default style. ---
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "290.31"
$c.ClearFormats()
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.3675"
$c.ClearFormats()
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "49.47"
$c.ClearFormats()
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.3389"
$c.ClearFormats()
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "1.170"
$c.ClearFormats()
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.07612"
$c.ClearFormats()
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "6.053"
$c.ClearFormats()
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "6.899"
$c.ClearFormats()
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.00001134"
$c.ClearFormats()
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "89.12"
$c.ClearFormats()
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "0.06764"
$c.ClearFormats()
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "6.241"
$c.ClearFormats()
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "0.5327"
$c.ClearFormats()
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "16.49"
$c.ClearFormats()
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "12.02"
$c.ClearFormats()
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "2.385"
$c.ClearFormats()
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "2.981"
$c.ClearFormats()
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "19.90"
$c.ClearFormats()
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "145.72"
$c.ClearFormats()
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "4.970"
$c.ClearFormats()
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "1.032"
$c.ClearFormats()
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "6.248"
$c.ClearFormats()
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "1.993"
$c.ClearFormats()
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "10.29"
$c.ClearFormats()
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.08457"
$c.ClearFormats()
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.02544"
$c.ClearFormats()
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.2327"
$c.ClearFormats()
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.06563"
$c.ClearFormats()
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "5.544"
$c.ClearFormats()
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.6375"
$c.ClearFormats()
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "14.38"
$c.ClearFormats()
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.6017"
$c.ClearFormats()
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "123.49"
$c.ClearFormats()
